$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, shifting existing rows 92:190 down to 93:191
$ws.Rows("92").Insert()

# Populate the newly inserted row 92 with the new data entry.
# Columns A,B,C,E,F,G,H,I,J,K,T hold constant boilerplate values identical
# to every other data row in this sheet.
$ws.Cells.Item(92, 1).Value = 4
$ws.Cells.Item(92, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(92, 3).Value = "Los Lagos"
$ws.Cells.Item(92, 4).Value = 44587
$ws.Cells.Item(92, 5).Value = 10
$ws.Cells.Item(92, 6).Value = "Fruta"
$ws.Cells.Item(92, 7).Value = 100101
$ws.Cells.Item(92, 8).Value = "Berries"
$ws.Cells.Item(92, 9).Value = 100112025
$ws.Cells.Item(92, 10).Value = "Frutilla"
$ws.Cells.Item(92, 11).Value = "Sin especificar"
$ws.Cells.Item(92, 12).Value = "Primera"
$ws.Cells.Item(92, 13).Value = 80
$ws.Cells.Item(92, 14).Value = 8500
$ws.Cells.Item(92, 15).Value = 9000
$ws.Cells.Item(92, 16).Value = 8750
$ws.Cells.Item(92, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(92, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(92, 19).Value = 1250
$ws.Cells.Item(92, 20).Value = 7
